$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 979.5
$ws.Range("I28").Value = 1030.7059
$ws.Range("J28").Value = 900.36365
$ws.Range("K28").Value = 1030.7059
$ws.Range("L28").Value = 900.36365
$ws.Range("M28").Value = -545.7058999999999
$ws.Range("N28").Value = -1870.36365

$ws.Range("H92").Value = 754.4666999999999
$ws.Range("I92").Value = 485.53845
$ws.Range("J92").Value = 2502.5
$ws.Range("K92").Value = 485.53845
$ws.Range("L92").Value = 2502.5
$ws.Range("M92").Value = 762.46155
$ws.Range("N92").Value = -4998.5

$ws.Range("H94").Value = 2362.4
$ws.Range("I94").Value = 2362.4
$ws.Range("K94").Value = 2362.4
$ws.Range("M94").Value = -1911.4

$ws.Range("H96").Value = 279.82352
$ws.Range("I96").Value = 276
$ws.Range("K96").Value = 828
$ws.Range("M96").Value = 545

$ws.Range("H97").Value = 1609.875
$ws.Range("J97").Value = 1899.8334
$ws.Range("L97").Value = 5699.5002
$ws.Range("N97").Value = -6691.5002

$ws.Range("H100").Value = 2084.3333
$ws.Range("I100").Value = 1900
$ws.Range("K100").Value = 1900
$ws.Range("M100").Value = -1359

$ws.Range("H101").Value = 1163.1875
$ws.Range("J101").Value = 1502.5
$ws.Range("L101").Value = 4507.5
$ws.Range("N101").Value = -7751.5

$ws.Range("H103").Value = 931.75
$ws.Range("I103").Value = 717.05554
$ws.Range("J103").Value = 1575.8334
$ws.Range("K103").Value = 2151.16662
$ws.Range("L103").Value = 4727.5002
$ws.Range("M103").Value = -1565.16662
$ws.Range("N103").Value = -5899.5002

$ws.Range("H132").Value = 846661.6
$ws.Range("I132").Value = 1717.102
$ws.Range("J132").Value = 5446915
$ws.Range("K132").Value = 5151.306
$ws.Range("L132").Value = 16340745
$ws.Range("M132").Value = -2621.306
$ws.Range("N132").Value = -16345805

$ws.Range("H135").Value = 19221.436
$ws.Range("I135").Value = 21885.979
$ws.Range("K135").Value = 196973.811
$ws.Range("M135").Value = -194438.811

$ws.Range("H137").Value = 3450415.2
$ws.Range("I137").Value = 12502380
$ws.Range("J137").Value = 2047.4762
$ws.Range("K137").Value = 37507140
$ws.Range("L137").Value = 6142.4286
$ws.Range("M137").Value = -37504590
$ws.Range("N137").Value = -11242.4286

$ws.Range("H138").Value = 2733987.2
$ws.Range("I138").Value = 1298.35
$ws.Range("J138").Value = 7939109
$ws.Range("K138").Value = 3895.05
$ws.Range("L138").Value = 23817327
$ws.Range("M138").Value = 1244.95
$ws.Range("N138").Value = -23827607

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1431.421
$ws.Range("I45").Value = 969.7692
$ws.Range("J45").Value = 2431.6667
$ws.Range("K45").Value = 969.7692
$ws.Range("L45").Value = 2431.6667
$ws.Range("M45").Value = -592.7692
$ws.Range("N45").Value = -3185.6667

$ws.Range("H61").Value = 77077770
$ws.Range("I61").Value = 100100696
$ws.Range("J61").Value = 334671.34
$ws.Range("K61").Value = 100100696
$ws.Range("L61").Value = 334671.34
$ws.Range("M61").Value = -100100484
$ws.Range("N61").Value = -335095.34

$ws.Range("H74").Value = 6001251.5
$ws.Range("I74").Value = 8097959
$ws.Range("J74").Value = 92350
$ws.Range("K74").Value = 8097959
$ws.Range("L74").Value = 92350
$ws.Range("M74").Value = -8097085
$ws.Range("N74").Value = -94098

$ws.Range("H77").Value = 6001251.5
$ws.Range("I77").Value = 8097959
$ws.Range("J77").Value = 92350
$ws.Range("K77").Value = 40489795
$ws.Range("L77").Value = 461750
$ws.Range("M77").Value = -40485427
$ws.Range("N77").Value = -470486

$ws.Range("H97").Value = 1563106.6
$ws.Range("I97").Value = 2155732.8
$ws.Range("J97").Value = 728.2727
$ws.Range("K97").Value = 2155732.8
$ws.Range("L97").Value = 728.2727
$ws.Range("M97").Value = -2155236.8
$ws.Range("N97").Value = -1720.2727

$ws.Range("H102").Value = 7525827.5
$ws.Range("I102").Value = 8936358
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 8936358
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -8934736
$ws.Range("N102").Value = -6244

$ws.Range("H110").Value = 1378.6666
$ws.Range("I110").Value = 1105
$ws.Range("J110").Value = 1743.5555
$ws.Range("K110").Value = 1105
$ws.Range("L110").Value = 1743.5555
$ws.Range("M110").Value = 940
$ws.Range("N110").Value = -5833.5555

$ws.Range("H123").Value = 50814.832
$ws.Range("J123").Value = 50814.832
$ws.Range("L123").Value = 50814.832
$ws.Range("N123").Value = -60614.832

$ws.Range("H132").Value = 59125.305
$ws.Range("I132").Value = 40888.48
$ws.Range("K132").Value = 122665.44
$ws.Range("M132").Value = -120135.44

$ws.Range("H136").Value = 77077770
$ws.Range("I136").Value = 100100696
$ws.Range("J136").Value = 334671.34
$ws.Range("K136").Value = 300302088
$ws.Range("L136").Value = 1004014.02
$ws.Range("M136").Value = -300299538
$ws.Range("N136").Value = -1009114.02

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 449.6111
$ws.Range("I94").Value = 292.35715
$ws.Range("K94").Value = 292.35715
$ws.Range("M94").Value = 158.64285

$ws.Range("H99").Value = 1468.5714
$ws.Range("I99").Value = 1413.3334
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 1413.3334
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = 84.66660000000002
$ws.Range("N99").Value = -4796

$ws.Range("H107").Value = 1850.8334
$ws.Range("I107").Value = 1818.2941
$ws.Range("J107").Value = 1893.3846
$ws.Range("K107").Value = 1818.2941
$ws.Range("L107").Value = 1893.3846
$ws.Range("M107").Value = 101.7058999999999
$ws.Range("N107").Value = -5733.3846

$ws.Range("H134").Value = 1871.2
$ws.Range("I134").Value = 1024.3334
$ws.Range("J134").Value = 2839.0476
$ws.Range("K134").Value = 3073.0002
$ws.Range("L134").Value = 8517.1428
$ws.Range("M134").Value = -538.0001999999999
$ws.Range("N134").Value = -13587.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2076.027
$ws.Range("I31").Value = 994.36536
$ws.Range("J31").Value = 4632.6816
$ws.Range("K31").Value = 994.36536
$ws.Range("L31").Value = 4632.6816
$ws.Range("M31").Value = -699.36536
$ws.Range("N31").Value = -5222.6816

$ws.Range("H34").Value = 2076.027
$ws.Range("I34").Value = 994.36536
$ws.Range("J34").Value = 4632.6816
$ws.Range("K34").Value = 994.36536
$ws.Range("L34").Value = 4632.6816
$ws.Range("M34").Value = -792.36536
$ws.Range("N34").Value = -5036.6816

$ws.Range("H52").Value = 49500
$ws.Range("J52").Value = 49500
$ws.Range("L52").Value = 49500
$ws.Range("N52").Value = -50088

$ws.Range("H58").Value = 40002784
$ws.Range("I58").Value = 52634140
$ws.Range("J58").Value = 3485.8333
$ws.Range("K58").Value = 52634140
$ws.Range("L58").Value = 3485.8333
$ws.Range("M58").Value = -52633937
$ws.Range("N58").Value = -3891.8333

$ws.Range("H105").Value = 1343.8334
$ws.Range("I105").Value = 1318.5333
$ws.Range("K105").Value = 1318.5333
$ws.Range("M105").Value = 428.4666999999999

$ws.Range("H107").Value = 765.5
$ws.Range("I107").Value = 622.2
$ws.Range("J107").Value = 1004.3333
$ws.Range("K107").Value = 622.2
$ws.Range("L107").Value = 1004.3333
$ws.Range("M107").Value = 1297.8
$ws.Range("N107").Value = -4844.3333

$ws.Range("H136").Value = 40002784
$ws.Range("I136").Value = 52634140
$ws.Range("J136").Value = 3485.8333
$ws.Range("K136").Value = 157902420
$ws.Range("L136").Value = 10457.4999
$ws.Range("M136").Value = -157899870
$ws.Range("N136").Value = -15557.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6832.727

$ws.Range("H52").Value = 33212.918
$ws.Range("I52").Value = 2930
$ws.Range("J52").Value = 35965.91
$ws.Range("K52").Value = 8790
$ws.Range("L52").Value = 107897.73
$ws.Range("M52").Value = -8524
$ws.Range("N52").Value = -108429.73

$ws.Range("H122").Value = 711.76666
$ws.Range("I122").Value = 271.625
$ws.Range("J122").Value = 1214.7858
$ws.Range("K122").Value = 2444.625
$ws.Range("L122").Value = 10933.0722
$ws.Range("M122").Value = 5.375
$ws.Range("N122").Value = -15833.0722

$ws.Range("H131").Value = 913.25
$ws.Range("I131").Value = 432.15384
$ws.Range("J131").Value = 1481.8182
$ws.Range("K131").Value = 1296.46152
$ws.Range("L131").Value = 4445.4546
$ws.Range("M131").Value = 3743.53848
$ws.Range("N131").Value = -14525.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2183.3333
$ws.Range("I97").Value = 2209.0908
$ws.Range("J97").Value = 1900
$ws.Range("K97").Value = 2209.0908
$ws.Range("L97").Value = 1900
$ws.Range("M97").Value = -1713.0908
$ws.Range("N97").Value = -2892

$ws.Range("H107").Value = 429.22223
$ws.Range("I107").Value = 206.66667
$ws.Range("J107").Value = 540.5
$ws.Range("K107").Value = 206.66667
$ws.Range("L107").Value = 540.5
$ws.Range("M107").Value = 1713.33333
$ws.Range("N107").Value = -4380.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 9800
$ws.Range("J42").Value = 9800
$ws.Range("L42").Value = 9800
$ws.Range("N42").Value = -10926

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("N43").Value = 0

$ws.Range("H49").Value = 9800
$ws.Range("J49").Value = 9800
$ws.Range("L49").Value = 9800
$ws.Range("N49").Value = -10094

$ws.Range("H93").Value = 400
$ws.Range("I93").Value = 400
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 400
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = 848

$ws.Range("H122").Value = 3290.6086
$ws.Range("I122").Value = 2494.889
$ws.Range("J122").Value = 3802.1428
$ws.Range("K122").Value = 7484.667
$ws.Range("L122").Value = 11406.4284
$ws.Range("M122").Value = -5034.667
$ws.Range("N122").Value = -16306.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 48496.523
$ws.Range("I100").Value = 45839.91
$ws.Range("J100").Value = 51418.8
$ws.Range("K100").Value = 91679.82000000001
$ws.Range("L100").Value = 102837.6
$ws.Range("M100").Value = -103919.6
$ws.Range("N100").Value = -103919.6
